{"js": "// 1) Update the arraignment date from July 06, 2022 to July 09, 2022.\nconst dateResults = context.document.body.search(\" on July 06, 2022.\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\n\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(\" on July 09, 2022.\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) Remove the \"...and report to jail on October 14, 2022, at 7:00 p.m\" clause,\n//    so the sentence ends with \"...by October 18, 2022.\" instead.\nconst startResults = context.document.body.search(\"October 18, 2022\", { matchCase: true });\nstartResults.load(\"items\");\nawait context.sync();\n\nconst endResults = context.document.body.search(\"p.m\", { matchCase: true });\nendResults.load(\"items\");\nawait context.sync();\n\nif (startResults.items.length > 0 && endResults.items.length > 0) {\n  const afterOctober18 = startResults.items[0].getRange(Word.RangeLocation.end);\n  const throughPm = endResults.items[0];\n  const clauseToRemove = afterOctober18.expandTo(throughPm);\n  clauseToRemove.delete();\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Update the arraignment date from July 06, 2022 to July 09, 2022.\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Text = \" on July 06, 2022.\"\n$find1.Replacement.ClearFormatting()\n$find1.Replacement.Text = \" on July 09, 2022.\"\n$find1.Execute([ref]\" on July 06, 2022.\", $false, $false, $false, $false, $false, $true, 1, $false, \" on July 09, 2022.\", 2)\n\n# 2) Remove the \"...and report to jail on October 14, 2022, at 7:00 p.m\" clause,\n#    so the sentence ends with \"...by October 18, 2022.\" instead.\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = \"October 18, 2022\"\n$found2 = $find2.Execute()\n$afterOctober18 = $find2.Parent.Duplicate\n$afterOctober18.Collapse(0)\n\n$find3 = $d.Content.Find\n$find3.ClearFormatting()\n$find3.Text = \"p.m\"\n$found3 = $find3.Execute()\n$throughPm = $find3.Parent\n\n$clauseToRemove = $d.Range($afterOctober18.Start, $throughPm.End)\n$clauseToRemove.Delete()\n"}
